$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

# Clear the old "Status"/"status" header cells in G1:I1 but keep their (blank) styles
$ws.Range("G1:I1").ClearContents()

# Clear the old "Pass"/"pass"/"text N" cells in G2:I5 entirely (no style, so they disappear)
$ws.Range("G2:I5").ClearContents()

# Write the new "status" header across J1:S1
$ws.Range("J1:S1").Value2 = "status"

# Write the new "pass" values down J2:J5
$ws.Range("J2:J5").Value2 = "pass"

# Give each of J1:S1 its own distinct (but visually default) cell style,
# mirroring how the workbook ends up with 10 extra near-duplicate style records
$cols = @("J","K","L","M","N","O","P","Q","R","S")
$i = 0
foreach ($c in $cols) {
  $i = $i + 1
  $cell = $ws.Range($c + "1")
  $cell.IndentLevel = $i
}

# Approximate the resulting best-fit column widths for the shifted/new columns
$ws.Columns.Item(7).ColumnWidth = 5.66
$ws.Columns.Item(8).ColumnWidth = 5.66
$ws.Columns.Item(9).ColumnWidth = 5.66
$ws.Columns.Item(10).ColumnWidth = 5.66

# Update the active selection to C6, as recorded in the saved view state
$ws.Activate()
[void]$ws.Range("C6").Select()
